$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cfh"
$ws.Cells.Item(2,3).Value = "Itgam"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.573304666666667
$ws.Cells.Item(2,8).Value = 4.719914
$ws.Cells.Item(2,9).Value = 0.02227938362376841
$ws.Cells.Item(2,10).Value = 0.0222793836237684
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 66.38494866666666
$ws.Cells.Item(2,14).Value = 199.154846
$ws.Cells.Item(2,15).Value = 0.4171200956172241
$ws.Cells.Item(2,16).Value = 0.4171200956172242
$ws.Cells.Item(2,17).Value = 104.4437495336938
$ws.Cells.Item(2,18).Value = 939.993745803244
$ws.Cells.Item(2,19).Value = 0.009293178627439096
$ws.Cells.Item(2,20).Value = 0.009293178627439096

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cfh"
$ws.Cells.Item(3,3).Value = "Itgam"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.573304666666667
$ws.Cells.Item(3,8).Value = 4.719914
$ws.Cells.Item(3,9).Value = 0.02227938362376841
$ws.Cells.Item(3,10).Value = 0.0222793836237684
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 92.76573566666667
$ws.Cells.Item(3,14).Value = 278.297207
$ws.Cells.Item(3,15).Value = 0.5828799043827758
$ws.Cells.Item(3,16).Value = 0.5828799043827758
$ws.Cells.Item(3,17).Value = 145.9487648311331
$ws.Cells.Item(3,18).Value = 1313.538883480198
$ws.Cells.Item(3,19).Value = 0.01298620499632931
$ws.Cells.Item(3,20).Value = 0.01298620499632931

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Cfh"
$ws.Cells.Item(4,3).Value = "Itgam"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 36.796554
$ws.Cells.Item(4,8).Value = 110.389662
$ws.Cells.Item(4,9).Value = 0.5210717033819111
$ws.Cells.Item(4,10).Value = 0.521071703381911
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 66.38494866666666
$ws.Cells.Item(4,14).Value = 199.154846
$ws.Cells.Item(4,15).Value = 0.4171200956172241
$ws.Cells.Item(4,16).Value = 0.4171200956172242
$ws.Cells.Item(4,17).Value = 2442.737348400228
$ws.Cells.Item(4,18).Value = 21984.63613560205
$ws.Cells.Item(4,19).Value = 0.2173494787380926
$ws.Cells.Item(4,20).Value = 0.2173494787380926

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Cfh"
$ws.Cells.Item(5,3).Value = "Itgam"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 36.796554
$ws.Cells.Item(5,8).Value = 110.389662
$ws.Cells.Item(5,9).Value = 0.5210717033819111
$ws.Cells.Item(5,10).Value = 0.521071703381911
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 92.76573566666667
$ws.Cells.Item(5,14).Value = 278.297207
$ws.Cells.Item(5,15).Value = 0.5828799043827758
$ws.Cells.Item(5,16).Value = 0.5828799043827758
$ws.Cells.Item(5,17).Value = 3413.459401808226
$ws.Cells.Item(5,18).Value = 30721.13461627404
$ws.Cells.Item(5,19).Value = 0.3037222246438185
$ws.Cells.Item(5,20).Value = 0.3037222246438184

# Row 6
$ws.Cells.Item(6,1).Value = "M2"
$ws.Cells.Item(6,2).Value = "Cfh"
$ws.Cells.Item(6,3).Value = "Itgam"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.68148033333333
$ws.Cells.Item(6,8).Value = 53.04444099999999
$ws.Cells.Item(6,9).Value = 0.2503853778156443
$ws.Cells.Item(6,10).Value = 0.2503853778156443
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 66.38494866666666
$ws.Cells.Item(6,14).Value = 199.154846
$ws.Cells.Item(6,15).Value = 0.4171200956172241
$ws.Cells.Item(6,16).Value = 0.4171200956172242
$ws.Cells.Item(6,17).Value = 1173.784164279009
$ws.Cells.Item(6,18).Value = 10564.05747851108
$ws.Cells.Item(6,19).Value = 0.1044407727356164
$ws.Cells.Item(6,20).Value = 0.1044407727356164

# Row 7
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Cfh"
$ws.Cells.Item(7,3).Value = "Itgam"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.68148033333333
$ws.Cells.Item(7,8).Value = 53.04444099999999
$ws.Cells.Item(7,9).Value = 0.2503853778156443
$ws.Cells.Item(7,10).Value = 0.2503853778156443
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 92.76573566666667
$ws.Cells.Item(7,14).Value = 278.297207
$ws.Cells.Item(7,15).Value = 0.5828799043827758
$ws.Cells.Item(7,16).Value = 0.5828799043827758
$ws.Cells.Item(7,17).Value = 1640.235530797365
$ws.Cells.Item(7,18).Value = 14762.11977717629
$ws.Cells.Item(7,19).Value = 0.145944605080028
$ws.Cells.Item(7,20).Value = 0.145944605080028

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Cfh"
$ws.Cells.Item(8,3).Value = "Itgam"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 14.56572533333333
$ws.Cells.Item(8,8).Value = 43.697176
$ws.Cells.Item(8,9).Value = 0.2062635351786761
$ws.Cells.Item(8,10).Value = 0.2062635351786761
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 66.38494866666666
$ws.Cells.Item(8,14).Value = 199.154846
$ws.Cells.Item(8,15).Value = 0.4171200956172241
$ws.Cells.Item(8,16).Value = 0.4171200956172242
$ws.Cells.Item(8,17).Value = 966.9449285460994
$ws.Cells.Item(8,18).Value = 8702.504356914895
$ws.Cells.Item(8,19).Value = 0.08603666551607604
$ws.Cells.Item(8,20).Value = 0.08603666551607604

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Cfh"
$ws.Cells.Item(9,3).Value = "Itgam"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 14.56572533333333
$ws.Cells.Item(9,8).Value = 43.697176
$ws.Cells.Item(9,9).Value = 0.2062635351786761
$ws.Cells.Item(9,10).Value = 0.2062635351786761
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 92.76573566666667
$ws.Cells.Item(9,14).Value = 278.297207
$ws.Cells.Item(9,15).Value = 0.5828799043827758
$ws.Cells.Item(9,16).Value = 0.5828799043827758
$ws.Cells.Item(9,17).Value = 1351.20022606527
$ws.Cells.Item(9,18).Value = 12160.80203458743
$ws.Cells.Item(9,19).Value = 0.1202268696626
$ws.Cells.Item(9,20).Value = 0.1202268696626
